# Insert a new weekly price record as the second data row (row 22) for the
# "Terminal La Palmera de La Serena - Perejil" sheet. Inserting the row
# shifts every following row down by one (old row 22 -> 23, ..., old row
# 119 -> 120), which matches the rest of the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 22..119 down to 23..120, leaving a blank row 22 behind.
$ws.Rows("22:22").Insert()

# Populate the new row 22 with the new weekly record.
$ws.Range("A22").Value = 8
$ws.Range("B22").Value = "Terminal La Palmera de La Serena"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44600
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 100112044
$ws.Range("G22").Value = "Perejil"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 2200
$ws.Range("K22").Value = 2300
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = 2400
$ws.Range("N22").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O22").Value = "Provincia del Elquí"
$ws.Range("P22").Value = 1600
$ws.Range("Q22").Value = 1.5
$ws.Range("R22").Value = "Hortaliza"
